$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = "Dr. Nahla Nagiub, Dr. Mohammad El-Tanany, Dr. Servinaz Sayed Mohammad, Dr. Nesma, Dr. Rana Abo-Zaid"
$ws.Range("G3").Value = "Dr. Hend Mahmoud, Dr. Amira Sobhy, Dr. Veronia Rafat, Dr. Servinaz Sayed Mohammad, Dr. Menna tuâ€™Allah Medhat, Dr. Rana Abo-Zaid, Dr. Asmaa Reda"
$ws.Range("G4").Value = "Dr. Hend Mahmoud, Dr. Amira Sobhy, Dr. Majorelle Magdy, Dr. Gehan Adel, Dr. Manar Montaser, Dr. Menna tuâ€™Allah Medhat, Dr. Rana Abo-Zaid, Dr. Asmaa Reda"
$ws.Range("G5").Value = "Dr. Lamiaa Ossama, Dr. Amera Ahmad Saad, Dr. Menna tu'Alllah Mohammad, Dr. Nada Gouda, Dr. Fatma Elhady"
$ws.Range("G6").Value = "Dr. Lamiaa Ossama, Dr. Kerelos Zareef, Dr. Amera Ahmad Saad, Dr. Abeer Ragab, Dr. Menna tu'Alllah Mohammad, Dr. Nada Mohammad, Dr. Fatma Elhady"
$ws.Range("G8").Value = "Dr. Dalia Mohammad Abd Al-Salam, Dr. Marwa Mustafa, Dr. Dina Adel, Dr. Amira Ibrahim, Dr. Madeha Saeed"
$ws.Range("G9").Value = "Dr. Merna Said, Dr. Yasmeena Fattoh, Dr. Maryam Ahmad, Dr. Madeha Saeed, Dr. Arwa Al-Sayed"
$ws.Range("G10").Value = "Dr. Al-Shimaa Khaled, Dr. Mohammad Safwat"
$ws.Range("G11").Value = "D Mariam E. Mohammad, Dr. Mona Ibrahim Hussein, Dr. Heba Al-Sayed Mohammad"
$ws.Range("G12").Value = "D Mariam E. Mohammad, Dr. Mona Ibrahim Hussein, Dr. Heba Al-Sayed Mohammad"
$ws.Range("G15").Value = "Dr. Amr Saeed, Dr. Walaa Ghanima"
$ws.Range("G16").Value = "Dr. Abdullah El-Agrody, Dr. Neveen Nashaat, Dr. Remon, Dr. Wafaa Ebida, Dr. Salma Hassan, Dr. Nardine, Dr. Aya Hanafy, Dr. Eman Samir Gabry"
$ws.Range("G17").Value = "Dr. Neveen Nashaat, Dr. Wafaa Ebida, Dr. Marina Sorial, Dr. Yasmin, Dr. Eman Samir Gabry"
$ws.Range("G18").Value = "Dr. Nahla Nagiub, Dr. Mohammad El-Tanany, Dr. Servinaz Sayed Mohammad, Dr. Nesma, Dr. Rana Abo-Zaid"
$ws.Range("G19").Value = "Dr. Amira Sobhy, Dr. Veronia Rafat, Dr. Mohammad El-Tanany, Dr. Eman Tantawi, Dr. Servinaz Sayed Mohammad, Dr. Asmaa Reda"
$ws.Range("G20").Value = "Dr. Amira Sobhy, Dr. Mohammad El-Tanany, Dr. Servinaz Sayed Mohammad, Dr. Nesma, Dr. Nourhan Mahmoud, Dr. Heba Mahmoud Ali, Dr. Asmaa Reda"
$ws.Range("G21").Value = "Dr. Lamiaa Ossama, Dr. Amera Ahmad Saad, Dr. Menna tu'Alllah Mohammad, Dr. Nada Gouda, Dr. Fatma Elhady"
$ws.Range("G22").Value = "Dr. Lamiaa Ossama, Dr. Kerelos Zareef, Dr. Amera Ahmad Saad, Dr. Abeer Ragab, Dr. Menna tu'Alllah Mohammad, Dr. Nada Mohammad, Dr. Fatma Elhady"
$ws.Range("G24").Value = "Dr. Dalia Mohammad Abd Al-Salam, Dr. Marwa Mustafa, Dr. Dina Adel, Dr. Amira Ibrahim, Dr. Madeha Saeed"
$ws.Range("G25").Value = "Dr. Merna Said, Dr. Yasmeena Fattoh, Dr. Maryam Ahmad, Dr. Madeha Saeed, Dr. Arwa Al-Sayed"
$ws.Range("G26").Value = "Dr. Al-Shimaa Khaled, Dr. Mohammad Safwat"
$ws.Range("G27").Value = "D Mariam E. Mohammad, Dr. Mona Ibrahim Hussein, Dr. Heba Al-Sayed Mohammad"
$ws.Range("G28").Value = "D Mariam E. Mohammad, Dr. Mona Ibrahim Hussein, Dr. Heba Al-Sayed Mohammad"
$ws.Range("G31").Value = "Dr. Amr Saeed, Dr. Walaa Ghanima"
$ws.Range("G32").Value = "Dr. Abdullah El-Agrody, Dr. Neveen Nashaat, Dr. Remon, Dr. Wafaa Ebida, Dr. Salma Hassan, Dr. Nardine, Dr. Aya Hanafy, Dr. Eman Samir Gabry"
$ws.Range("G33").Value = "Dr. Neveen Nashaat, Dr. Wafaa Ebida, Dr. Marina Sorial, Dr. Yasmin, Dr. Eman Samir Gabry"
$ws.Range("G34").Value = "Dr. Amira Sobhy, Dr. Nahla Nagiub, Dr. Veronia Rafat, Dr. Gehan Adel, Dr. Servinaz Sayed Mohammad, Dr. Menna tuâ€™Allah Medhat, Dr. Rana Abo-Zaid, Dr. Asmaa Reda, Administrator"
$ws.Range("G35").Value = "Dr. Mohammad El-Tanany, Dr. Eman Tantawi, Dr. Gehan Adel, Dr. Servinaz Sayed Mohammad, Dr. Rana Abo-Zaid, Dr. Heba Mahmoud Ali, Administrator"
$ws.Range("G36").Value = "Dr. Hend Mahmoud, Dr. Alshimaa Atef, Dr. Veronia Rafat, Dr. Amira Sobhy, Dr. Servinaz Sayed Mohammad, Dr. Shimaa Ahmad Mekki, Dr. Rana Abo-Zaid, Dr. Heba Mahmoud Ali"
$ws.Range("G37").Value = "Dr. Lamiaa Ossama, Dr. Kerelos Zareef, Dr. Abeer Ragab, Dr. Nada Mohammad, Dr. Nada Gouda, Dr. Fatma Elhady"
$ws.Range("G40").Value = "Dr. Merna Said, Dr. Sara Atawia, Dr. Eman M. Abo-Sakaya, Dr. Mai Mustafa, Dr. Yasmeena Fattoh, Dr. Nahed Mosaad, Dr. Nourhan Osama, Dr. Basma Hamed, Dr. Maryam Ahmad, Dr. Madeha Saeed, Dr. Amany Raafat, Dr. Merna Mahrous, Dr. Marina Youhanna"
$ws.Range("G41").Value = "Dr. Eman M. Abo-Sakaya, Dr. Dina Adel, Dr. Esraa Mostafa, Dr. Nadia Mostafa, Dr. Nourhan Osama, Dr. Maryam Ahmad, Dr. Amira Ibrahim, Dr. Sarah Abdelmohsen, Dr. Amany Raafat, Dr. Merna Mahrous"
$ws.Range("G43").Value = "Dr. Youstina Gamil, Dr. Mona Ibrahim Hussein, Dr. Sarah Mahdy"
$ws.Range("G44").Value = "Dr. Youstina Gamil, Dr. Mona Ibrahim Hussein, Dr. Sarah Mahdy"
$ws.Range("G46").Value = "Dr. Nourham Mostafa, Dr. Aya Alaa-Eldein"
$ws.Range("G48").Value = "Dr. Maryam Ashraf, Dr. Remon"
$ws.Range("G49").Value = "Dr. Neveen Nashaat, Dr. Ola Abd Al-Fattah, Dr. Eman Samir Gabry, Dr. Remon, Dr. Monica, Dr. Naema Gomaa"
$ws.Range("G50").Value = "Dr. Amira Sobhy, Dr. Nahla Nagiub, Dr. Veronia Rafat, Dr. Gehan Adel, Dr. Servinaz Sayed Mohammad, Dr. Menna tuâ€™Allah Medhat, Dr. Rana Abo-Zaid, Dr. Asmaa Reda, Administrator"
$ws.Range("G51").Value = "Dr. Mohammad El-Tanany, Dr. Eman Tantawi, Dr. Gehan Adel, Dr. Servinaz Sayed Mohammad, Dr. Rana Abo-Zaid, Dr. Heba Mahmoud Ali, Administrator"
$ws.Range("G52").Value = "Dr. Hend Mahmoud, Dr. Alshimaa Atef, Dr. Veronia Rafat, Dr. Amira Sobhy, Dr. Servinaz Sayed Mohammad, Dr. Shimaa Ahmad Mekki, Dr. Rana Abo-Zaid, Dr. Heba Mahmoud Ali"
$ws.Range("G53").Value = "Dr. Lamiaa Ossama, Dr. Kerelos Zareef, Dr. Abeer Ragab, Dr. Nada Mohammad, Dr. Nada Gouda, Dr. Fatma Elhady"
$ws.Range("G56").Value = "Dr. Merna Said, Dr. Sara Atawia, Dr. Eman M. Abo-Sakaya, Dr. Mai Mustafa, Dr. Yasmeena Fattoh, Dr. Nahed Mosaad, Dr. Nourhan Osama, Dr. Basma Hamed, Dr. Maryam Ahmad, Dr. Madeha Saeed, Dr. Amany Raafat, Dr. Merna Mahrous, Dr. Marina Youhanna"
$ws.Range("G57").Value = "Dr. Eman M. Abo-Sakaya, Dr. Dina Adel, Dr. Esraa Mostafa, Dr. Nadia Mostafa, Dr. Nourhan Osama, Dr. Maryam Ahmad, Dr. Amira Ibrahim, Dr. Sarah Abdelmohsen, Dr. Amany Raafat, Dr. Merna Mahrous"
$ws.Range("G59").Value = "Dr. Youstina Gamil, Dr. Mona Ibrahim Hussein, Dr. Sarah Mahdy"
$ws.Range("G60").Value = "Dr. Youstina Gamil, Dr. Mona Ibrahim Hussein, Dr. Sarah Mahdy"
$ws.Range("G62").Value = "Dr. Nourham Mostafa, Dr. Aya Alaa-Eldein"
$ws.Range("G64").Value = "Dr. Maryam Ashraf, Dr. Remon"
$ws.Range("G65").Value = "Dr. Neveen Nashaat, Dr. Ola Abd Al-Fattah, Dr. Eman Samir Gabry, Dr. Remon, Dr. Monica, Dr. Naema Gomaa"
$ws.Range("G66").Value = "Dr. Hend Mahmoud, Dr. Amira Sobhy, Dr. Veronia Rafat, Dr. Nahla Nagiub, Dr. Eman Tantawi, Dr. Gehan Adel, Dr. Servinaz Sayed Mohammad, Dr. Menna tuâ€™Allah Medhat, Dr. Nourhan Mahmoud, Dr. Asmaa Reda"
$ws.Range("G67").Value = "Dr. Hend Mahmoud, Dr. Amira Sobhy, Dr. Veronia Rafat, Dr. Servinaz Sayed Mohammad, Dr. Menna tuâ€™Allah Medhat, Dr. Rana Abo-Zaid, Dr. Asmaa Reda"
$ws.Range("G68").Value = "Dr. Alshimaa Atef, Dr. Veronia Rafat, Dr. Amira Sobhy, Dr. Eman Tantawi, Dr. Nourhan Mahmoud"
$ws.Range("G69").Value = "Dr. Lamiaa Ossama, Dr. Amera Ahmad Saad, Dr. Abeer Ragab, Dr. Menna tu'Alllah Mohammad, Dr. Nada Mohammad"
$ws.Range("G70").Value = "Dr. Amera Ahmad Saad, Dr. Fatma Elhady, Dr. Nada Gouda"
$ws.Range("G72").Value = "Dr. Merna Said, Dr. Sara Atawia, Dr. Eman M. Abo-Sakaya, Dr. Mai Mustafa, Dr. Yasmeena Fattoh, Dr. Nahed Mosaad, Dr. Nourhan Osama, Dr. Basma Hamed, Dr. Maryam Ahmad, Dr. Madeha Saeed, Dr. Amany Raafat, Dr. Merna Mahrous, Dr. Marina Youhanna"
$ws.Range("G73").Value = "Dr. Merna Said, Dr. Dalia Mohammad Abd Al-Salam, Dr. Esraa Mostafa, Dr. Nahed Mosaad, Dr. Madeha Saeed, Dr. Arwa Al-Sayed"
$ws.Range("G74").Value = "Dr. Al-Shimaa Khaled, Dr. Mohammad Safwat"
$ws.Range("G79").Value = "Dr. Amr Saeed, Dr. Walaa Ghanima"
$ws.Range("G80").Value = "Dr. Eman Mohammad Al, Dr. Neveen Nashaat, Dr. Ola Abd Al-Fattah, Dr. Salma Hassan, Dr. Aya Hanafy, Dr. Marina Atef"
$ws.Range("G81").Value = "Dr. Neveen Nashaat, Dr. Wafaa Ebida, Dr. Marina Sorial, Dr. Yasmin, Dr. Eman Samir Gabry"
$ws.Range("G82").Value = "Dr. Hend Mahmoud, Dr. Amira Sobhy, Dr. Veronia Rafat, Dr. Nahla Nagiub, Dr. Eman Tantawi, Dr. Gehan Adel, Dr. Servinaz Sayed Mohammad, Dr. Menna tuâ€™Allah Medhat, Dr. Nourhan Mahmoud, Dr. Asmaa Reda"
$ws.Range("G83").Value = "Dr. Amira Sobhy, Dr. Veronia Rafat, Dr. Mohammad El-Tanany, Dr. Eman Tantawi, Dr. Servinaz Sayed Mohammad, Dr. Asmaa Reda"
$ws.Range("G84").Value = "Dr. Amira Sobhy, Dr. Mohammad El-Tanany, Dr. Servinaz Sayed Mohammad, Dr. Nesma, Dr. Nourhan Mahmoud, Dr. Heba Mahmoud Ali, Dr. Asmaa Reda"
$ws.Range("G85").Value = "Dr. Lamiaa Ossama, Dr. Amera Ahmad Saad, Dr. Abeer Ragab, Dr. Menna tu'Alllah Mohammad, Dr. Nada Mohammad"
$ws.Range("G86").Value = "Dr. Amera Ahmad Saad, Dr. Fatma Elhady, Dr. Nada Gouda"
$ws.Range("G88").Value = "Dr. Merna Said, Dr. Sara Atawia, Dr. Eman M. Abo-Sakaya, Dr. Mai Mustafa, Dr. Yasmeena Fattoh, Dr. Nahed Mosaad, Dr. Nourhan Osama, Dr. Basma Hamed, Dr. Maryam Ahmad, Dr. Madeha Saeed, Dr. Amany Raafat, Dr. Merna Mahrous, Dr. Marina Youhanna"
$ws.Range("G89").Value = "Dr. Merna Said, Dr. Dalia Mohammad Abd Al-Salam, Dr. Esraa Mostafa, Dr. Nahed Mosaad, Dr. Madeha Saeed, Dr. Arwa Al-Sayed"
$ws.Range("G90").Value = "Dr. Al-Shimaa Khaled, Dr. Mohammad Safwat"
$ws.Range("G95").Value = "Dr. Amr Saeed, Dr. Walaa Ghanima"
$ws.Range("G96").Value = "Dr. Eman Mohammad Al, Dr. Neveen Nashaat, Dr. Ola Abd Al-Fattah, Dr. Salma Hassan, Dr. Aya Hanafy, Dr. Marina Atef"
$ws.Range("G98").Value = "Dr. Nahla Nagiub, Dr. Mohammad El-Tanany, Dr. Servinaz Sayed Mohammad, Dr. Nesma, Dr. Rana Abo-Zaid"
$ws.Range("G100").Value = "Dr. Alshimaa Atef, Dr. Veronia Rafat, Dr. Amira Sobhy, Dr. Eman Tantawi, Dr. Nourhan Mahmoud"
$ws.Range("G101").Value = "Dr. Lamiaa Ossama, Dr. Kerelos Zareef, Dr. Abeer Ragab, Dr. Nada Mohammad, Dr. Nada Gouda, Dr. Fatma Elhady"
$ws.Range("G102").Value = "Dr. Amera Ahmad Saad, Dr. Fatma Elhady, Dr. Nada Gouda"
$ws.Range("G104").Value = "Dr. Eman M. Abo-Sakaya, Dr. Dina Adel, Dr. Esraa Mostafa, Dr. Yasmeena Fattoh, Dr. Nourhan Osama, Dr. Maryam Ahmad, Dr. Amira Ibrahim, Dr. Arwa Al-Sayed, Dr. Marina Youhanna"
$ws.Range("G105").Value = "Dr. Eman M. Abo-Sakaya, Dr. Dina Adel, Dr. Esraa Mostafa, Dr. Nadia Mostafa, Dr. Nourhan Osama, Dr. Maryam Ahmad, Dr. Amira Ibrahim, Dr. Sarah Abdelmohsen, Dr. Amany Raafat, Dr. Merna Mahrous"
$ws.Range("G107").Value = "Dr. Youstina Gamil, Dr. Mona Ibrahim Hussein, Dr. Sarah Mahdy"
$ws.Range("G108").Value = "Dr. Youstina Gamil, Dr. Mona Ibrahim Hussein, Dr. Sarah Mahdy"
$ws.Range("G112").Value = "Dr. Yassmen Ahmad, Dr. Neveen Nashaat, Dr. Nahla, Dr. Youstina Magdy, Dr. Salma Hassan, Dr. Remon"
$ws.Range("G114").Value = "Dr. Nahla Nagiub, Dr. Mohammad El-Tanany, Dr. Servinaz Sayed Mohammad, Dr. Nesma, Dr. Rana Abo-Zaid"
$ws.Range("G116").Value = "Dr. Hend Mahmoud, Dr. Amira Sobhy, Dr. Majorelle Magdy, Dr. Gehan Adel, Dr. Manar Montaser, Dr. Menna tuâ€™Allah Medhat, Dr. Rana Abo-Zaid, Dr. Asmaa Reda"
$ws.Range("G117").Value = "Dr. Lamiaa Ossama, Dr. Kerelos Zareef, Dr. Abeer Ragab, Dr. Nada Mohammad, Dr. Nada Gouda, Dr. Fatma Elhady"
$ws.Range("G118").Value = "Dr. Amera Ahmad Saad, Dr. Fatma Elhady, Dr. Nada Gouda"
$ws.Range("G120").Value = "Dr. Eman M. Abo-Sakaya, Dr. Dina Adel, Dr. Esraa Mostafa, Dr. Yasmeena Fattoh, Dr. Nourhan Osama, Dr. Maryam Ahmad, Dr. Amira Ibrahim, Dr. Arwa Al-Sayed, Dr. Marina Youhanna"
$ws.Range("G121").Value = "Dr. Eman M. Abo-Sakaya, Dr. Dina Adel, Dr. Esraa Mostafa, Dr. Nadia Mostafa, Dr. Nourhan Osama, Dr. Maryam Ahmad, Dr. Amira Ibrahim, Dr. Sarah Abdelmohsen, Dr. Amany Raafat, Dr. Merna Mahrous"
$ws.Range("G123").Value = "Dr. Youstina Gamil, Dr. Mona Ibrahim Hussein, Dr. Sarah Mahdy"
$ws.Range("G124").Value = "Dr. Youstina Gamil, Dr. Mona Ibrahim Hussein, Dr. Sarah Mahdy"
$ws.Range("G128").Value = "Dr. Yassmen Ahmad, Dr. Neveen Nashaat, Dr. Nahla, Dr. Youstina Magdy, Dr. Salma Hassan, Dr. Remon"
